$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the two "Check response..." hyperlinks (the new layout has no
#    hyperlinks at all).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 2. Insert a new column before the old "Test Case" column (B) to make
#    room for the new "Test Case Number" column. This shifts the old
#    B/C/D columns to C/D/E and - conveniently - the inserted column
#    picks up the same row-banding style (s=1/s=2) already used by the
#    rest of each row, so no extra style plumbing is required for it.
# ---------------------------------------------------------------------
$ws.Columns("B:B").Insert()

# ---------------------------------------------------------------------
# 3. Header row
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Test Case Number"

# ---------------------------------------------------------------------
# 4. New "Test Case Number" values (0-based running index per row)
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4

# ---------------------------------------------------------------------
# 5. Normalise the Email/password columns (D/E) for every data row so
#    they all reference the same sample e-mail address / phone number,
#    matching rows 2 and 4 which already had the right data.
# ---------------------------------------------------------------------
$ws.Range("D3").Value = "wiasm.mtour@gmail.com"
$ws.Range("E3").Value = 123456789

$ws.Range("D5").Value = "wiasm.mtour@gmail.com"
$ws.Range("E5").Value = 123456789

$ws.Range("D6").Value = "wiasm.mtour@gmail.com"
$ws.Range("E6").Value = 123456789

# ---------------------------------------------------------------------
# 6. Re-apply the correct (pre-existing) cell style to the cells whose
#    content we just replaced, by copying formats only from a cell that
#    already carries the desired style - this reuses the existing style
#    index instead of minting new ones in styles.xml.
# ---------------------------------------------------------------------
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)

$ws.Range("E3").Copy()
$ws.Range("E6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 7. Column widths: widen the "Test Case" column and size the new
#    password column.
# ---------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 67.83333333333333
$ws.Columns("E").ColumnWidth = 28.833333333333332

# ---------------------------------------------------------------------
# 8. Selection / view state
# ---------------------------------------------------------------------
$ws.Range("D5:E6").Select()

Write-Output "done"
